# commit: "commit hawaii samples, 30 from 20190910, rechecking CRM and last sample"
#
# Data entry of a new CRM-accuracy sample (2019-09-10 / serial 43718) into
# row 48 of Sheet1, plus updating the window/selection state to reflect
# where the user was working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook window was minimized when this was saved.
$excel.ActiveWindow.WindowState = -4140   # xlMinimized

# Row 48 was blank except for the CRM reference value (C48) and the
# shared "% off" formula (D48). Fill in the new sample's date, batch
# value, sample size and batch note, matching the formatting already
# used for the rows above (A47 carries the m/d/yyyy date style).
$ws.Range("A47").Copy()
$ws.Range("A48").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A48").Value = 43718                 # 9/10/2019
$ws.Range("B48").Value = 2218.7914332390801    # batch value
$ws.Range("E48").Value = 169                   # sample size
$ws.Range("F48").Value = $ws.Range("F47").Value2   # "Opened CRM (9/8/2019), Dani B"

# Leave the selection where the user last clicked.
$ws.Range("E50").Select()
